$wb = $excel.ActiveWorkbook

# Row 40 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6168.1665
$ws.Range("J40").Value = 6502.8
$ws.Range("L40").Value = 6502.8
$ws.Range("N40").Value = -6852.8

# Row 41 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 973.8333

# Row 62 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6668.6665
$ws.Range("I62").Value = 9000
$ws.Range("K62").Value = 9000
$ws.Range("M62").Value = -8376

# Row 65 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 6668.6665
$ws.Range("I65").Value = 9000
$ws.Range("K65").Value = 45000
$ws.Range("M65").Value = -41880

# Row 125 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 54793.844
$ws.Range("J125").Value = 202497.6
$ws.Range("L125").Value = 1822478.4
$ws.Range("N125").Value = -1827398.4

# Row 137 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1624.7142

# Row 138 (ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2113.3542
$ws.Range("I138").Value = 2162.3333
$ws.Range("J138").Value = 2091.0908
$ws.Range("K138").Value = 6486.999899999999
$ws.Range("L138").Value = 6273.2724
$ws.Range("M138").Value = -1346.999899999999
$ws.Range("N138").Value = -16553.2724

# Row 32 (ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6041.8687
$ws.Range("I32").Value = 2378.9434
$ws.Range("K32").Value = 2378.9434
$ws.Range("M32").Value = -2091.9434

# Row 45 (ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2158.875
$ws.Range("I45").Value = 1442.8125
$ws.Range("K45").Value = 1442.8125
$ws.Range("M45").Value = -1065.8125

# Row 88 (ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 4320.4
$ws.Range("I88").Value = 3895
$ws.Range("J88").Value = 4426.75
$ws.Range("K88").Value = 3895
$ws.Range("L88").Value = 4426.75
$ws.Range("M88").Value = -3489
$ws.Range("N88").Value = -5238.75

# Row 91 (ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 4320.4
$ws.Range("I91").Value = 3895
$ws.Range("J91").Value = 4426.75
$ws.Range("K91").Value = 3895
$ws.Range("L91").Value = 4426.75
$ws.Range("M91").Value = -2491
$ws.Range("N91").Value = -7234.75

# Row 20 (BSM)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3122.5
$ws.Range("J20").Value = 4067.3333
$ws.Range("L20").Value = 4067.3333
$ws.Range("N20").Value = -4561.3333

# Row 105 (BSM)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3241.8
$ws.Range("I105").Value = 3177.375
$ws.Range("K105").Value = 3177.375
$ws.Range("M105").Value = -1430.375

# Row 31 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11159.75
$ws.Range("I31").Value = 4033
$ws.Range("K31").Value = 4033
$ws.Range("M31").Value = -3738

# Row 34 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 11159.75
$ws.Range("I34").Value = 4033
$ws.Range("K34").Value = 4033
$ws.Range("M34").Value = -3831

# Row 52 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 76633
$ws.Range("I52").Value = 71950
$ws.Range("J52").Value = 85999
$ws.Range("K52").Value = 71950
$ws.Range("L52").Value = 85999
$ws.Range("M52").Value = -71656
$ws.Range("N52").Value = -86587

# Row 58 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3212.2942
$ws.Range("I58").Value = 2619
$ws.Range("K58").Value = 2619
$ws.Range("M58").Value = -2416

# Row 86 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3594.5
$ws.Range("I86").Value = 3476
$ws.Range("J86").Value = 3950
$ws.Range("K86").Value = 3476
$ws.Range("L86").Value = 3950
$ws.Range("M86").Value = -2353
$ws.Range("N86").Value = -6196

# Row 89 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 3594.5
$ws.Range("I89").Value = 3476
$ws.Range("J89").Value = 3950
$ws.Range("K89").Value = 17380
$ws.Range("L89").Value = 19750
$ws.Range("M89").Value = -11764
$ws.Range("N89").Value = -30982

# Row 92 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 60000
$ws.Range("J92").Value = 60000
$ws.Range("L92").Value = 60000
$ws.Range("N92").Value = -64992

# Row 135 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 93897.5
$ws.Range("J135").Value = 93897.5
$ws.Range("L135").Value = 93897.5
$ws.Range("N135").Value = -104037.5

# Row 136 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3212.2942
$ws.Range("I136").Value = 2619
$ws.Range("K136").Value = 7857
$ws.Range("M136").Value = -5307

# Row 140 (CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 79943.25
$ws.Range("J140").Value = 79943.25
$ws.Range("L140").Value = 79943.25
$ws.Range("N140").Value = -90303.25

# Row 58 (CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 16730
$ws.Range("I58").Value = 16730
$ws.Range("K58").Value = 50190
$ws.Range("M58").Value = -50062

# Row 68 (CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2691
$ws.Range("I68").Value = 2702.8572
$ws.Range("K68").Value = 8108.571599999999
$ws.Range("M68").Value = -7297.571599999999

# Row 71 (CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2691
$ws.Range("I71").Value = 2702.8572
$ws.Range("K71").Value = 24325.7148
$ws.Range("M71").Value = -20269.7148

# Row 80 (GSM)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4472.143
$ws.Range("J80").Value = 4601.5
$ws.Range("L80").Value = 4601.5
$ws.Range("N80").Value = -6597.5

# Row 83 (GSM)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4472.143
$ws.Range("J83").Value = 4601.5
$ws.Range("L83").Value = 23007.5
$ws.Range("N83").Value = -32991.5

# Row 7 (LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6469.7646
$ws.Range("I7").Value = 5576.3076
$ws.Range("J7").Value = 9373.5
$ws.Range("K7").Value = 5576.3076
$ws.Range("L7").Value = 9373.5
$ws.Range("M7").Value = -5464.3076
$ws.Range("N7").Value = -9597.5

# Row 46 (LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3097.0908
$ws.Range("I46").Value = 1261
$ws.Range("K46").Value = 1261
$ws.Range("M46").Value = -1073

# Row 55 (LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 700.25
$ws.Range("I55").Value = 226
$ws.Range("J55").Value = 1581
$ws.Range("K55").Value = 226
$ws.Range("L55").Value = 1581
$ws.Range("M55").Value = -53
$ws.Range("N55").Value = -1927

# Row 74 (LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# Row 77 (LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# Row 126 (LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 6469.7646
$ws.Range("I126").Value = 5576.3076
$ws.Range("J126").Value = 9373.5
$ws.Range("K126").Value = 16728.9228
$ws.Range("L126").Value = 28120.5
$ws.Range("M126").Value = -14258.9228
$ws.Range("N126").Value = -33060.5

# Row 46 (WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 75107.25
$ws.Range("J46").Value = 75107.25
$ws.Range("L46").Value = 75107.25
$ws.Range("N46").Value = -75569.25

# Row 74 (WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 8187.75
$ws.Range("J74").Value = 10084
$ws.Range("L74").Value = 10084
$ws.Range("N74").Value = -11956

# Row 77 (WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H77").Value = 8187.75
$ws.Range("J77").Value = 10084
$ws.Range("L77").Value = 30252
$ws.Range("N77").Value = -39612

# Row 107 (WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1774.2413
$ws.Range("I107").Value = 1966.8334
$ws.Range("K107").Value = 5900.5002
$ws.Range("M107").Value = -3980.5002

# Row 134 (WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 75107.25
$ws.Range("J134").Value = 75107.25
$ws.Range("L134").Value = 225321.75
$ws.Range("N134").Value = -230391.75

# Row 136 (WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5668.091
$ws.Range("I136").Value = 4405.316
$ws.Range("K136").Value = 13215.948
$ws.Range("M136").Value = -10665.948

